$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: repeat of the header row (bold, centered), but only for
#     columns A, B, C and H (matches the "type"/"genome"/"URL" header
#     that precedes the new URL block) ---
$ws.Range("A9").Value = "type"
$ws.Range("B9").Value = "genome"
$ws.Range("C9").Value = "URL"
$ws.Range("H9").Value = "[this line is not part of the file format itself]"

# Pick up the bold/centered formatting from the existing header row (row 1)
# for each of the touched cells, so the shared style is reused instead of
# a brand-new style being created.
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("H9").PasteSpecial(-4122)

# --- Row 10: Barley URL row ---
$ws.Range("A10").Value = "URL"
$ws.Range("B10").Value = "Barley"
$ws.Range("C10").Value = "http://penguin.scri.ac.uk/paul/germinate/germinate_development/app/flapjack/flapjack_search/search.pl?marker="

# --- Row 11: Rice URL row ---
$ws.Range("A11").Value = "URL"
$ws.Range("B11").Value = "Rice"
$ws.Range("C11").Value = "http://rice.plantbiology.msu.edu/cgi-bin/gbrowse/rice/?name="

# --- Update the active selection to mirror the author's saved selection ---
[void]$ws.Range("H9").Select()
